$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix bug with duplicates: update the computed count/value in B3 from 9911 to 9905
$ws.Range("B3").Value = 9905
